# Updates market-price derived figures (columns H-N) on each job profit sheet,
# reflecting refreshed Universalis market data pulled by the scheduled runner.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 404.9091
$ws.Range("I4").Value = 295.4
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 295.4
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -181.4
$ws.Range("N4").Value = -1728
$ws.Range("H17").Value = 931.9583
$ws.Range("I17").Value = 199
$ws.Range("J17").Value = 998.5909
$ws.Range("K17").Value = 597
$ws.Range("L17").Value = 2995.7727
$ws.Range("M17").Value = -429
$ws.Range("N17").Value = -3331.7727
$ws.Range("H28").Value = 46993.1
$ws.Range("I28").Value = 453.4
$ws.Range("J28").Value = 93532.8
$ws.Range("K28").Value = 453.4
$ws.Range("L28").Value = 93532.8
$ws.Range("M28").Value = 31.60000000000002
$ws.Range("N28").Value = -94502.8
$ws.Range("H32").Value = 253356.5
$ws.Range("J32").Value = 253356.5
$ws.Range("L32").Value = 253356.5
$ws.Range("N32").Value = -254008.5
$ws.Range("H62").Value = 5430.8335
$ws.Range("I62").Value = 8001
$ws.Range("J62").Value = 4916.8
$ws.Range("K62").Value = 8001
$ws.Range("L62").Value = 4916.8
$ws.Range("M62").Value = -7377
$ws.Range("N62").Value = -6164.8
$ws.Range("H65").Value = 5430.8335
$ws.Range("I65").Value = 8001
$ws.Range("J65").Value = 4916.8
$ws.Range("K65").Value = 40005
$ws.Range("L65").Value = 24584
$ws.Range("M65").Value = -36885
$ws.Range("N65").Value = -30824
$ws.Range("H98").Value = 1993.125
$ws.Range("I98").Value = 1156.5
$ws.Range("K98").Value = 1156.5
$ws.Range("M98").Value = 341.5
$ws.Range("H112").Value = 995.6667
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 995.6667
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 2987.0001
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -5203.0001
$ws.Range("H122").Value = 1993.125
$ws.Range("I122").Value = 1156.5
$ws.Range("K122").Value = 3469.5
$ws.Range("M122").Value = -1019.5
$ws.Range("H132").Value = 2842.6316
$ws.Range("I132").Value = 2000.8235
$ws.Range("K132").Value = 6002.470499999999
$ws.Range("M132").Value = -3472.470499999999
$ws.Range("H138").Value = 2018
$ws.Range("I138").Value = 1669.0834
$ws.Range("J138").Value = 2483.2222
$ws.Range("K138").Value = 5007.2502
$ws.Range("L138").Value = 7449.6666
$ws.Range("M138").Value = 132.7497999999996
$ws.Range("N138").Value = -17729.6666
$ws.Range("H141").Value = 4562.077
$ws.Range("I141").Value = 3754
$ws.Range("K141").Value = 11262
$ws.Range("M141").Value = -6082

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 546.3333
$ws.Range("I2").Value = 424
$ws.Range("J2").Value = 974.5
$ws.Range("K2").Value = 424
$ws.Range("L2").Value = 974.5
$ws.Range("M2").Value = -311
$ws.Range("N2").Value = -1200.5
$ws.Range("H6").Value = 507499.75
$ws.Range("J6").Value = 1000000
$ws.Range("L6").Value = 1000000
$ws.Range("N6").Value = -1000346
$ws.Range("H32").Value = 5400.3804
$ws.Range("I32").Value = 2704.4807
$ws.Range("K32").Value = 2704.4807
$ws.Range("M32").Value = -2417.4807
$ws.Range("H45").Value = 7355435.5
$ws.Range("I45").Value = 2457.5715
$ws.Range("K45").Value = 2457.5715
$ws.Range("M45").Value = -2080.5715
$ws.Range("H97").Value = 573.6111
$ws.Range("I97").Value = 578.36365
$ws.Range("J97").Value = 566.1429000000001
$ws.Range("K97").Value = 578.36365
$ws.Range("L97").Value = 566.1429000000001
$ws.Range("M97").Value = -82.36365000000001
$ws.Range("N97").Value = -1558.1429
$ws.Range("H116").Value = 546.3333
$ws.Range("I116").Value = 424
$ws.Range("J116").Value = 974.5
$ws.Range("K116").Value = 424
$ws.Range("L116").Value = 974.5
$ws.Range("M116").Value = 1870
$ws.Range("N116").Value = -5562.5
$ws.Range("H122").Value = 3291.3333
$ws.Range("I122").Value = 3312
$ws.Range("K122").Value = 9936
$ws.Range("M122").Value = -7486

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 546.3333
$ws.Range("I3").Value = 424
$ws.Range("J3").Value = 974.5
$ws.Range("K3").Value = 424
$ws.Range("L3").Value = 974.5
$ws.Range("M3").Value = -310
$ws.Range("N3").Value = -1202.5
$ws.Range("H134").Value = 3044.0977
$ws.Range("I134").Value = 2386.7026
$ws.Range("K134").Value = 7160.1078
$ws.Range("M134").Value = -4625.1078

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 258.35294
$ws.Range("I7").Value = 254.66667
$ws.Range("J7").Value = 260.36365
$ws.Range("K7").Value = 254.66667
$ws.Range("L7").Value = 260.36365
$ws.Range("M7").Value = -141.66667
$ws.Range("N7").Value = -486.36365
$ws.Range("H86").Value = 6610.5557
$ws.Range("I86").Value = 4039
$ws.Range("K86").Value = 4039
$ws.Range("M86").Value = -2916
$ws.Range("H89").Value = 6610.5557
$ws.Range("I89").Value = 4039
$ws.Range("K89").Value = 20195
$ws.Range("M89").Value = -14579
$ws.Range("H99").Value = 7494941
$ws.Range("I99").Value = 12347875
$ws.Range("J99").Value = 3127299.5
$ws.Range("K99").Value = 12347875
$ws.Range("L99").Value = 3127299.5
$ws.Range("M99").Value = -12346377
$ws.Range("N99").Value = -3130295.5
$ws.Range("H126").Value = 7494941
$ws.Range("I126").Value = 12347875
$ws.Range("J126").Value = 3127299.5
$ws.Range("K126").Value = 37043625
$ws.Range("L126").Value = 9381898.5
$ws.Range("M126").Value = -37041155
$ws.Range("N126").Value = -9386838.5
$ws.Range("H132").Value = 2615
$ws.Range("I132").Value = 2367.889
$ws.Range("J132").Value = 3059.8
$ws.Range("K132").Value = 7103.667
$ws.Range("L132").Value = 9179.400000000001
$ws.Range("M132").Value = -4573.667
$ws.Range("N132").Value = -14239.4
$ws.Range("H134").Value = 3279.4814
$ws.Range("J134").Value = 3109.75
$ws.Range("L134").Value = 9329.25
$ws.Range("N134").Value = -14399.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 466.66666
$ws.Range("I2").Value = 466.66666
$ws.Range("K2").Value = 2799.99996
$ws.Range("M2").Value = -2686.99996
$ws.Range("H4").Value = 6366551.5
$ws.Range("I4").Value = 5263557
$ws.Range("K4").Value = 15790671
$ws.Range("M4").Value = -15790559
$ws.Range("H22").Value = 1821.4286
$ws.Range("I22").Value = 1125
$ws.Range("J22").Value = 2100
$ws.Range("K22").Value = 3375
$ws.Range("L22").Value = 6300
$ws.Range("M22").Value = -3206
$ws.Range("N22").Value = -6638
$ws.Range("H27").Value = 1821.4286
$ws.Range("I27").Value = 1125
$ws.Range("J27").Value = 2100
$ws.Range("K27").Value = 3375
$ws.Range("L27").Value = 6300
$ws.Range("M27").Value = -3273
$ws.Range("N27").Value = -6504
$ws.Range("H80").Value = 1799.4
$ws.Range("J80").Value = 1848.75
$ws.Range("L80").Value = 5546.25
$ws.Range("N80").Value = -7418.25
$ws.Range("H83").Value = 1799.4
$ws.Range("J83").Value = 1848.75
$ws.Range("L83").Value = 16638.75
$ws.Range("N83").Value = -25998.75
$ws.Range("H113").Value = 101004.3
$ws.Range("I113").Value = 850
$ws.Range("J113").Value = 112132.555
$ws.Range("K113").Value = 2550
$ws.Range("L113").Value = 336397.665
$ws.Range("M113").Value = -380
$ws.Range("N113").Value = -340737.665

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 244.375
$ws.Range("I2").Value = 23.333334
$ws.Range("J2").Value = 377
$ws.Range("K2").Value = 23.333334
$ws.Range("L2").Value = 377
$ws.Range("M2").Value = 89.66666599999999
$ws.Range("N2").Value = -603
$ws.Range("H102").Value = 1998.3077
$ws.Range("I102").Value = 1870.3
$ws.Range("J102").Value = 2425
$ws.Range("K102").Value = 1870.3
$ws.Range("L102").Value = 2425
$ws.Range("M102").Value = -248.3
$ws.Range("N102").Value = -5669
$ws.Range("H132").Value = 3313.0881
$ws.Range("I132").Value = 2608.818
$ws.Range("K132").Value = 7826.454000000001
$ws.Range("M132").Value = -5296.454000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3270858.8
$ws.Range("I40").Value = 2864.923
$ws.Range("J40").Value = 13891839
$ws.Range("K40").Value = 2864.923
$ws.Range("L40").Value = 13891839
$ws.Range("M40").Value = -2728.923
$ws.Range("N40").Value = -13892111
$ws.Range("H61").Value = 771
$ws.Range("I61").Value = 771
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 771
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -569
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 771
$ws.Range("I113").Value = 771
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 771
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1399
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 1441.9667
$ws.Range("I132").Value = 1320.963
$ws.Range("J132").Value = 2531
$ws.Range("K132").Value = 3962.889
$ws.Range("L132").Value = 7593
$ws.Range("M132").Value = -1432.889
$ws.Range("N132").Value = -12653
$ws.Range("H133").Value = 49818.285
$ws.Range("J133").Value = 49818.285
$ws.Range("L133").Value = 49818.285
$ws.Range("N133").Value = -54878.285
$ws.Range("H136").Value = 4247
$ws.Range("I136").Value = 5776.5
$ws.Range("K136").Value = 17329.5
$ws.Range("M136").Value = -14779.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 79999.664
$ws.Range("J108").Value = 79999.664
$ws.Range("L108").Value = 79999.664
$ws.Range("N108").Value = -87679.664
$ws.Range("H136").Value = 850.75
$ws.Range("I136").Value = 508
$ws.Range("K136").Value = 1524
$ws.Range("M136").Value = 1026

